$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (row 1) ---
# A1 "Customer Name" is unchanged; B1 changes from "Metadata" to "Line Detail"
$ws.Range("B1").Value = "Line Detail"

# --- Row 2: Company 1 (text already correct; JSON detail replaced) ---
$ws.Range("A2").Value = "Company 1"
$ws.Range("B2").Value = '[ {"Id": "1", "Desc": "Bolts", "Amount": 101.15}, {"Id": "2", "Desc": "Smith", "Amount": 90.80} ]'

# --- Row 3: Company 2 ---
$ws.Range("A3").Value = "Company 2"
$ws.Range("B3").Value = '[ {"Id": "1", "Desc": "Braces", "Amount": 51.15}, {"Id": "2", "Desc": "Wood", "Amount": 190.10} ]'

# --- Row 4: Company 3 (new) ---
$ws.Range("A4").Value = "Company 3"
$ws.Range("B4").Value = '[{"Id": "1", "Desc": "Braces", "Amount": 51.15}]'

# --- Row 5: Company 4 (new, no line-detail value) ---
$ws.Range("A5").Value = "Company 4"

# --- Row 6: Company 5 (new) ---
$ws.Range("A6").Value = "Company 5"
$ws.Range("B6").Value = '{"Id": "1", "Desc": "Braces", "Amount": 51.15}'

# B2 already carries the Menlo/9pt/gray "code" style from the source file;
# reuse it (via copy/paste-format) for the other detail cells instead of
# re-deriving the font, so the shared style table stays minimal/reused
# rather than growing new near-duplicate font entries.
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column widths (best-effort match of the authored widths; the host
#     quantizes ColumnWidth onto a coarse internal grid, so these are the
#     closest achievable inputs to the authored stored widths) ---
$ws.Columns.Item(1).ColumnWidth = 49.25
$ws.Columns.Item(2).ColumnWidth = 11.75
$ws.Columns.Item(3).ColumnWidth = 10.92
$ws.Columns.Item(4).ColumnWidth = 13.25
$ws.Columns.Item(5).ColumnWidth = 16.09

# --- Selection ends on B6, matching the authored state ---
[void]$ws.Range("B6").Select()
